# Applies the ATLAS tactics workbook update for release 4.7.0:
#  - Refresh the x-mitre-tactic STIX identifiers (column B) with new UUIDs
#  - Minor copy edits to a handful of tactic descriptions (column D)
#  - Bump the "created"/"last modified" date (columns F and G) to 29 October 2024

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tactics")

# --- New x-mitre-tactic STIX IDs for rows 2-15 (column B) ---
$newStixIds = @(
    "x-mitre-tactic--76d8d1bc-98fb-406a-a88d-a70649ad9365",
    "x-mitre-tactic--b7c3883a-1ae7-47cd-bd53-024c79ac833a",
    "x-mitre-tactic--35a4b685-8a9f-4f66-bef0-92c8a5fd8411",
    "x-mitre-tactic--b3f5fda4-a31a-4a74-9588-995b143d3436",
    "x-mitre-tactic--b777d4af-fbe8-4e44-9917-c5c5542a7147",
    "x-mitre-tactic--55e80a13-be94-43ba-aa57-e6c5cb822864",
    "x-mitre-tactic--4386734e-fa57-4700-95b8-e76b8cab8ab3",
    "x-mitre-tactic--14ef1e2d-f902-450a-a7e8-b032e318bb86",
    "x-mitre-tactic--528bc2da-8855-44c0-9e45-ab92c179bf6f",
    "x-mitre-tactic--2666e534-bf47-4656-9404-e90afe41f4a7",
    "x-mitre-tactic--4a088a36-6786-4486-a3ea-3c3576d61daa",
    "x-mitre-tactic--ef06a48d-1ccb-42c6-b5dd-2770a58f02d8",
    "x-mitre-tactic--9ceaa8fe-57f1-4923-a4a5-121b5111139c",
    "x-mitre-tactic--1fa998a1-d720-488c-b156-b569fe4e6308"
)

for ($i = 0; $i -lt $newStixIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newStixIds[$i]
}

# --- Description text fixes (column D) ---

# Row 6: Execution - link "Remote System Discovery"
$ws.Range("D6").Value = "The adversary is trying to run malicious code embedded in machine learning artifacts or software.`n`nExecution consists of techniques that result in adversary-controlled code running on a local or remote system.`nTechniques that run malicious code are often paired with techniques from all other tactics to achieve broader goals, like exploring a network or stealing data.`nFor example, an adversary might use a remote access tool to run a PowerShell script that does [Remote System Discovery](https://attack.mitre.org/techniques/T1018/).`n"

# Row 10: ML Attack Staging - "manor" -> "manner"
$ws.Range("D10").Value = "The adversary is leveraging their knowledge of and access to the target system to tailor the attack.`n`nML Attack Staging consists of techniques adversaries use to prepare their attack on the target ML model.`nTechniques can include training proxy models, poisoning the target model, and crafting adversarial data to feed the target model.`nSome of these techniques can be performed in an offline manner and are thus difficult to mitigate.`nThese techniques are often used to achieve the adversary's end goal.`n"

# Row 14: Reconnaissance - "organizations" -> "organizations'"
$ws.Range("D14").Value = "The adversary is trying to gather information about the machine learning system they can use to plan future operations.`n`nReconnaissance consists of techniques that involve adversaries actively or passively gathering information that can be used to support targeting.`nSuch information may include details of the victim organizations' machine learning capabilities and research efforts.`nThis information can be leveraged by the adversary to aid in other phases of the adversary lifecycle, such as using gathered information to obtain relevant ML artifacts, targeting ML capabilities used by the victim, tailoring attacks to the particular models used by the victim, or to drive and lead further Reconnaissance efforts.`n"

# Row 15: Resource Development - link "ML Attack Staging"
$ws.Range("D15").Value = "The adversary is trying to establish resources they can use to support operations.`n`nResource Development consists of techniques that involve adversaries creating,`npurchasing, or compromising/stealing resources that can be used to support targeting.`nSuch resources include machine learning artifacts, infrastructure, accounts, or capabilities.`nThese resources can be leveraged by the adversary to aid in other phases of the adversary lifecycle, such as [ML Attack Staging](/tactics/AML.TA0001).`n"

# --- Update created / last modified date columns (F and G) for all data rows ---
$newDate = "29 October 2024"
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 6).Value = $newDate
    $ws.Cells.Item($row, 7).Value = $newDate
}
